$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 480.04166
$ws.Cells.Item(53, 9).Value = 324.7143
$ws.Cells.Item(53, 10).Value = 697.5
$ws.Cells.Item(53, 11).Value = 324.7143
$ws.Cells.Item(53, 12).Value = 697.5
$ws.Cells.Item(53, 13).Value = 312.2857
$ws.Cells.Item(53, 14).Value = -1971.5

$ws.Cells.Item(132, 8).Value = 4172.48
$ws.Cells.Item(132, 9).Value = 4157.0435
$ws.Cells.Item(132, 11).Value = 12471.1305
$ws.Cells.Item(132, 13).Value = -9941.130499999999

$ws.Cells.Item(137, 8).Value = 6817
$ws.Cells.Item(137, 9).Value = 5509.4
$ws.Cells.Item(137, 10).Value = 7822.846
$ws.Cells.Item(137, 11).Value = 16528.2
$ws.Cells.Item(137, 12).Value = 23468.538
$ws.Cells.Item(137, 13).Value = -13978.2
$ws.Cells.Item(137, 14).Value = -28568.538

$ws.Cells.Item(138, 8).Value = 43480036
$ws.Cells.Item(138, 9).Value = 1286.6923
$ws.Cells.Item(138, 10).Value = 100002420
$ws.Cells.Item(138, 11).Value = 3860.0769
$ws.Cells.Item(138, 12).Value = 300007260
$ws.Cells.Item(138, 13).Value = 1279.9231
$ws.Cells.Item(138, 14).Value = -300017540

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10756201
$ws.Cells.Item(32, 9).Value = 12823739
$ws.Cells.Item(32, 10).Value = 5003.067
$ws.Cells.Item(32, 11).Value = 12823739
$ws.Cells.Item(32, 12).Value = 5003.067
$ws.Cells.Item(32, 13).Value = -12823452
$ws.Cells.Item(32, 14).Value = -5577.067

$ws.Cells.Item(61, 8).Value = 40004420
$ws.Cells.Item(61, 9).Value = 55558564
$ws.Cells.Item(61, 10).Value = 8054.4287
$ws.Cells.Item(61, 11).Value = 55558564
$ws.Cells.Item(61, 12).Value = 8054.4287
$ws.Cells.Item(61, 13).Value = -55558352
$ws.Cells.Item(61, 14).Value = -8478.4287

$ws.Cells.Item(74, 8).Value = 41716396
$ws.Cells.Item(74, 9).Value = 45508496
$ws.Cells.Item(74, 10).Value = 3295
$ws.Cells.Item(74, 11).Value = 45508496
$ws.Cells.Item(74, 12).Value = 3295
$ws.Cells.Item(74, 13).Value = -45507622
$ws.Cells.Item(74, 14).Value = -5043

$ws.Cells.Item(77, 8).Value = 41716396
$ws.Cells.Item(77, 9).Value = 45508496
$ws.Cells.Item(77, 10).Value = 3295
$ws.Cells.Item(77, 11).Value = 227542480
$ws.Cells.Item(77, 12).Value = 16475
$ws.Cells.Item(77, 13).Value = -227538112
$ws.Cells.Item(77, 14).Value = -25211

$ws.Cells.Item(132, 8).Value = 26381134
$ws.Cells.Item(132, 9).Value = 1791.6207
$ws.Cells.Item(132, 10).Value = 111381240
$ws.Cells.Item(132, 11).Value = 5374.8621
$ws.Cells.Item(132, 12).Value = 334143720
$ws.Cells.Item(132, 13).Value = -2844.8621
$ws.Cells.Item(132, 14).Value = -334148780

$ws.Cells.Item(136, 8).Value = 40004420
$ws.Cells.Item(136, 9).Value = 55558564
$ws.Cells.Item(136, 10).Value = 8054.4287
$ws.Cells.Item(136, 11).Value = 166675692
$ws.Cells.Item(136, 12).Value = 24163.2861
$ws.Cells.Item(136, 13).Value = -166673142
$ws.Cells.Item(136, 14).Value = -29263.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 20022.363
$ws.Cells.Item(86, 9).Value = 12606
$ws.Cells.Item(86, 10).Value = 39799.332
$ws.Cells.Item(86, 11).Value = 12606
$ws.Cells.Item(86, 12).Value = 39799.332
$ws.Cells.Item(86, 13).Value = -11483
$ws.Cells.Item(86, 14).Value = -42045.332

$ws.Cells.Item(89, 8).Value = 20022.363
$ws.Cells.Item(89, 9).Value = 12606
$ws.Cells.Item(89, 10).Value = 39799.332
$ws.Cells.Item(89, 11).Value = 63030
$ws.Cells.Item(89, 12).Value = 198996.66
$ws.Cells.Item(89, 13).Value = -57414
$ws.Cells.Item(89, 14).Value = -210228.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 24395784
$ws.Cells.Item(31, 9).Value = 3673.875
$ws.Cells.Item(31, 10).Value = 40006732
$ws.Cells.Item(31, 11).Value = 3673.875
$ws.Cells.Item(31, 12).Value = 40006732
$ws.Cells.Item(31, 13).Value = -3378.875
$ws.Cells.Item(31, 14).Value = -40007322

$ws.Cells.Item(34, 8).Value = 24395784
$ws.Cells.Item(34, 9).Value = 3673.875
$ws.Cells.Item(34, 10).Value = 40006732
$ws.Cells.Item(34, 11).Value = 3673.875
$ws.Cells.Item(34, 12).Value = 40006732
$ws.Cells.Item(34, 13).Value = -3471.875
$ws.Cells.Item(34, 14).Value = -40007136

$ws.Cells.Item(58, 8).Value = 2962.6072
$ws.Cells.Item(58, 9).Value = 2304.6
$ws.Cells.Item(58, 10).Value = 4607.625
$ws.Cells.Item(58, 11).Value = 2304.6
$ws.Cells.Item(58, 12).Value = 4607.625
$ws.Cells.Item(58, 13).Value = -2101.6
$ws.Cells.Item(58, 14).Value = -5013.625

$ws.Cells.Item(107, 8).Value = 660.8
$ws.Cells.Item(107, 9).Value = 152.5
$ws.Cells.Item(107, 11).Value = 152.5
$ws.Cells.Item(107, 13).Value = 1767.5

$ws.Cells.Item(132, 8).Value = 2277.125
$ws.Cells.Item(132, 9).Value = 2323.9565
$ws.Cells.Item(132, 10).Value = 1200
$ws.Cells.Item(132, 11).Value = 6971.869499999999
$ws.Cells.Item(132, 12).Value = 3600
$ws.Cells.Item(132, 13).Value = -4441.869499999999
$ws.Cells.Item(132, 14).Value = -8660

$ws.Cells.Item(136, 8).Value = 2962.6072
$ws.Cells.Item(136, 9).Value = 2304.6
$ws.Cells.Item(136, 10).Value = 4607.625
$ws.Cells.Item(136, 11).Value = 6913.799999999999
$ws.Cells.Item(136, 12).Value = 13822.875
$ws.Cells.Item(136, 13).Value = -4363.799999999999
$ws.Cells.Item(136, 14).Value = -18922.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(36, 8).Value = 2134.8235
$ws.Cells.Item(36, 9).Value = 566
$ws.Cells.Item(36, 10).Value = 3899.75
$ws.Cells.Item(36, 11).Value = 1698
$ws.Cells.Item(36, 12).Value = 11699.25
$ws.Cells.Item(36, 13).Value = -1529
$ws.Cells.Item(36, 14).Value = -12037.25

$ws.Cells.Item(56, 8).Value = 57334.168
$ws.Cells.Item(56, 9).Value = 57334.168
$ws.Cells.Item(56, 11).Value = 57334.168
$ws.Cells.Item(56, 13).Value = -56804.168

$ws.Cells.Item(113, 8).Value = 1974.6842
$ws.Cells.Item(113, 9).Value = 509.75
$ws.Cells.Item(113, 10).Value = 2365.3333
$ws.Cells.Item(113, 11).Value = 1529.25
$ws.Cells.Item(113, 12).Value = 7095.999899999999
$ws.Cells.Item(113, 13).Value = 640.75
$ws.Cells.Item(113, 14).Value = -11435.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1375.6666
$ws.Cells.Item(80, 9).Value = 1137
$ws.Cells.Item(80, 10).Value = 1853
$ws.Cells.Item(80, 11).Value = 1137
$ws.Cells.Item(80, 12).Value = 1853
$ws.Cells.Item(80, 13).Value = -139
$ws.Cells.Item(80, 14).Value = -3849

$ws.Cells.Item(83, 8).Value = 1375.6666
$ws.Cells.Item(83, 9).Value = 1137
$ws.Cells.Item(83, 10).Value = 1853
$ws.Cells.Item(83, 11).Value = 5685
$ws.Cells.Item(83, 12).Value = 9265
$ws.Cells.Item(83, 13).Value = -693
$ws.Cells.Item(83, 14).Value = -19249

$ws.Cells.Item(97, 8).Value = 2161.8262
$ws.Cells.Item(97, 9).Value = 2024.909
$ws.Cells.Item(97, 10).Value = 2287.3333
$ws.Cells.Item(97, 11).Value = 2024.909
$ws.Cells.Item(97, 12).Value = 2287.3333
$ws.Cells.Item(97, 13).Value = -1528.909
$ws.Cells.Item(97, 14).Value = -3279.3333

$ws.Cells.Item(113, 8).Value = 3537.9614
$ws.Cells.Item(113, 9).Value = 2777.5
$ws.Cells.Item(113, 10).Value = 4425.1665
$ws.Cells.Item(113, 11).Value = 2777.5
$ws.Cells.Item(113, 12).Value = 4425.1665
$ws.Cells.Item(113, 13).Value = -607.5
$ws.Cells.Item(113, 14).Value = -8765.166499999999

$ws.Cells.Item(132, 8).Value = 1541.5
$ws.Cells.Item(132, 9).Value = 1481.174
$ws.Cells.Item(132, 11).Value = 4443.522
$ws.Cells.Item(132, 13).Value = -1913.522

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3304.9473
$ws.Cells.Item(16, 9).Value = 3340.8235
$ws.Cells.Item(16, 11).Value = 3340.8235
$ws.Cells.Item(16, 13).Value = -3170.8235

$ws.Cells.Item(22, 8).Value = 1198.825
$ws.Cells.Item(22, 9).Value = 1091.3529
$ws.Cells.Item(22, 10).Value = 1278.2609
$ws.Cells.Item(22, 11).Value = 1091.3529
$ws.Cells.Item(22, 12).Value = 1278.2609
$ws.Cells.Item(22, 13).Value = -796.3529000000001
$ws.Cells.Item(22, 14).Value = -1868.2609

$ws.Cells.Item(27, 8).Value = 1198.825
$ws.Cells.Item(27, 9).Value = 1091.3529
$ws.Cells.Item(27, 10).Value = 1278.2609
$ws.Cells.Item(27, 11).Value = 1091.3529
$ws.Cells.Item(27, 12).Value = 1278.2609
$ws.Cells.Item(27, 13).Value = -984.3529000000001
$ws.Cells.Item(27, 14).Value = -1492.2609

$ws.Cells.Item(55, 8).Value = 584.2222
$ws.Cells.Item(55, 9).Value = 256.22223
$ws.Cells.Item(55, 11).Value = 256.22223
$ws.Cells.Item(55, 13).Value = -83.22223000000002

$ws.Cells.Item(64, 8).Value = 54287.5
$ws.Cells.Item(64, 10).Value = 54287.5
$ws.Cells.Item(64, 12).Value = 54287.5
$ws.Cells.Item(64, 14).Value = -54737.5

$ws.Cells.Item(67, 8).Value = 54287.5
$ws.Cells.Item(67, 10).Value = 54287.5
$ws.Cells.Item(67, 12).Value = 54287.5
$ws.Cells.Item(67, 14).Value = -55847.5

$ws.Cells.Item(93, 8).Value = 2238.6667
$ws.Cells.Item(93, 9).Value = 2049.8
$ws.Cells.Item(93, 10).Value = 2474.75
$ws.Cells.Item(93, 11).Value = 2049.8
$ws.Cells.Item(93, 12).Value = 2474.75
$ws.Cells.Item(93, 13).Value = -801.8000000000002
$ws.Cells.Item(93, 14).Value = -4970.75

$ws.Cells.Item(136, 8).Value = 2340.6216
$ws.Cells.Item(136, 9).Value = 2337.9429
$ws.Cells.Item(136, 10).Value = 2387.5
$ws.Cells.Item(136, 11).Value = 7013.8287
$ws.Cells.Item(136, 12).Value = 7162.5
$ws.Cells.Item(136, 13).Value = -4463.8287
$ws.Cells.Item(136, 14).Value = -12262.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6990.5
$ws.Cells.Item(62, 10).Value = 7092.5557
$ws.Cells.Item(62, 12).Value = 7092.5557
$ws.Cells.Item(62, 14).Value = -8340.555700000001

$ws.Cells.Item(65, 8).Value = 6990.5
$ws.Cells.Item(65, 10).Value = 7092.5557
$ws.Cells.Item(65, 12).Value = 35462.7785
$ws.Cells.Item(65, 14).Value = -41702.7785

$ws.Cells.Item(122, 8).Value = 55613236
$ws.Cells.Item(122, 9).Value = 91001144
$ws.Cells.Item(122, 10).Value = 3676.1428
$ws.Cells.Item(122, 11).Value = 273003432
$ws.Cells.Item(122, 12).Value = 11028.4284
$ws.Cells.Item(122, 13).Value = -273000982
$ws.Cells.Item(122, 14).Value = -15928.4284

$ws.Cells.Item(132, 8).Value = 1945.5476
$ws.Cells.Item(132, 9).Value = 1826.0322
$ws.Cells.Item(132, 10).Value = 2282.3635
$ws.Cells.Item(132, 11).Value = 5478.096600000001
$ws.Cells.Item(132, 12).Value = 6847.0905
$ws.Cells.Item(132, 13).Value = -2948.096600000001
$ws.Cells.Item(132, 14).Value = -11907.0905

$ws.Cells.Item(136, 8).Value = 3827.8572
$ws.Cells.Item(136, 9).Value = 1037.4615
$ws.Cells.Item(136, 10).Value = 8362.25
$ws.Cells.Item(136, 11).Value = 3112.3845
$ws.Cells.Item(136, 12).Value = 25086.75
$ws.Cells.Item(136, 13).Value = -562.3844999999997
$ws.Cells.Item(136, 14).Value = -30186.75
